$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values in this sheet are stored as text (inline strings) in the
# original workbook, even ones that look numeric (e.g. "238.11", "1.92").
# Force Text number format before assigning so Excel does not silently
# convert numeric-looking strings to actual numbers (which would also
# drop formatting like trailing zeros, e.g. "660.20" -> 660.2).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "97.276.75"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.706.09"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "238.11"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "1.92"
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("D7").Value = "660.20"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "3.702.72"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "0.0000310"
$ws.Range("E12").Value = "  +15.07%  "
$ws.Range("D13").Value = "44.34"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "6.79"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "4.393.16"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "97.057.35"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "9.18"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("D19").Value = "3.671.73"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "13.05"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("D23").Value = "523.33"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "3.44"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "101.84"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  +15.19%  "
$ws.Range("D29").Value = "13.64"
$ws.Range("E29").Value = "  +4.31%  "
$ws.Range("D30").Value = "3.902.00"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "12.69"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "3.07"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "32.30"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "653.90"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").Value = "0.598"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").Value = "8.89"
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("D43").Value = "2.07"
$ws.Range("E43").Value = "  +4.61%  "
$ws.Range("D44").Value = "6.83"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.491"
$ws.Range("E45").Value = "  +8.07%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "40.60"
$ws.Range("E46").Value = "  -8.96%  "
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "0.0463"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "23.63"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  +1.18%  "
